$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values are not
# auto-converted to the Number type (the source data stores them as strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.970.43"
$ws.Range("E2").Value = "  +1.91%  "

$ws.Range("D3").Value = "3.172.10"
$ws.Range("E3").Value = "  +4.05%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "579.63"
$ws.Range("E5").Value = "  +3.69%  "

$ws.Range("D6").Value = "151.67"
$ws.Range("E6").Value = "  +6.62%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.171.37"
$ws.Range("E8").Value = "  +4.11%  "

$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +3.40%  "

$ws.Range("E10").Value = "  +5.96%  "

$ws.Range("D11").Value = "6.25"
$ws.Range("E11").Value = "  -0.78%  "

$ws.Range("D12").Value = "0.502"
$ws.Range("E12").Value = "  +2.73%  "

$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  +18.02%  "

$ws.Range("D14").Value = "37.68"
$ws.Range("E14").Value = "  +5.69%  "

$ws.Range("D15").Value = "3.695.93"
$ws.Range("E15").Value = "  +4.16%  "

$ws.Range("D16").Value = "65.077.81"
$ws.Range("E16").Value = "  +1.97%  "

$ws.Range("D17").Value = "3.175.07"
$ws.Range("E17").Value = "  +3.96%  "

$ws.Range("D18").Value = "7.17"
$ws.Range("E18").Value = "  +5.57%  "

$ws.Range("E19").Value = "  +1.38%  "

$ws.Range("D20").Value = "514.06"
$ws.Range("E20").Value = "  +8.19%  "

$ws.Range("D21").Value = "14.85"
$ws.Range("E21").Value = "  +5.91%  "

$ws.Range("D22").Value = "0.730"
$ws.Range("E22").Value = "  +7.08%  "

$ws.Range("D23").Value = "15.27"
$ws.Range("E23").Value = "  +3.89%  "

$ws.Range("D24").Value = "7.80"

$ws.Range("D25").Value = "85.42"
$ws.Range("E25").Value = "  +3.15%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "8.97"
$ws.Range("E27").Value = "  +10.51%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "2.93"
$ws.Range("E28").Value = "  +4.89%  "

$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  +7.36%  "

$ws.Range("D30").Value = "27.89"
$ws.Range("E30").Value = "  +6.37%  "

$ws.Range("E31").Value = "  +13.89%  "

$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("E33").Value = "  +5.40%  "

$ws.Range("D34").Value = "6.32"
$ws.Range("E34").Value = "  +9.66%  "

$ws.Range("D35").Value = "6.57"
$ws.Range("E35").Value = "  +5.80%  "

$ws.Range("D36").Value = "55.75"

$ws.Range("D37").Value = "0.0896"
$ws.Range("E37").Value = "  +10.24%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "474.03"
$ws.Range("E38").Value = "  +6.51%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "3.14"
$ws.Range("E39").Value = "  +12.08%  "

$ws.Range("D40").Value = "0.0421"
$ws.Range("E40").Value = "  +2.78%  "

$ws.Range("D41").Value = "8.65"
$ws.Range("E41").Value = "  +4.57%  "

$ws.Range("D42").Value = "3.067.24"
$ws.Range("E42").Value = "  +1.47%  "

$ws.Range("D43").Value = "0.119"
$ws.Range("E43").Value = "  +1.32%  "

$ws.Range("E44").Value = "  +6.28%  "

$ws.Range("E45").Value = "  +6.26%  "

$ws.Range("D46").Value = "29.05"
$ws.Range("E46").Value = "  +3.06%  "

$ws.Range("D47").Value = "0.0₃0609"
$ws.Range("E47").Value = "  +18.96%  "

$ws.Range("E49").Value = "  +2.12%  "

$ws.Range("D50").Value = "2.25"
$ws.Range("E50").Value = "  +8.21%  "

$ws.Range("D51").Value = "120.54"
$ws.Range("E51").Value = "  +2.05%  "

# Restore the columns style so no stray number-format style is left
# referenced differently than the original workbook.
$ws.Range("D2:D51").Style = "Normal"
